# Auto-generated from xml_diff: update cryptos worksheet values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.934.67'
$ws.Range('E2').Value = '  +8.28%  '
$ws.Range('D3').Value = '2.558.77'
$ws.Range('E3').Value = '  +8.68%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '509.84'
$ws.Range('E5').Value = '  +6.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.80'
$ws.Range('E6').Value = '  +8.14%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.992'
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -2.94%  '
$ws.Range('D9').Value = '2.617.48'
$ws.Range('E9').Value = '  +11.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.15'
$ws.Range('E10').Value = '  +12.67%  '
$ws.Range('E11').Value = '  +7.69%  '
$ws.Range('E12').Value = '  +5.89%  '
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('D14').Value = '3.008.98'
$ws.Range('E14').Value = '  +8.93%  '
$ws.Range('D15').Value = '59.753.20'
$ws.Range('E15').Value = '  +8.19%  '
$ws.Range('E16').Value = '  +9.93%  '
$ws.Range('E17').Value = '  +7.27%  '
$ws.Range('D18').Value = '2.593.53'
$ws.Range('E18').Value = '  +10.29%  '
$ws.Range('E19').Value = '  +5.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '342.94'
$ws.Range('E20').Value = '  +8.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.51'
$ws.Range('E21').Value = '  +9.66%  '
$ws.Range('E22').Value = '  +8.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.20'
$ws.Range('E24').Value = '  +6.00%  '
$ws.Range('E25').Value = '  +6.89%  '
$ws.Range('E26').Value = '  +10.72%  '
$ws.Range('D27').Value = '2.679.61'
$ws.Range('E27').Value = '  +9.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.991'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('D29').Value = '0.0₃0843'
$ws.Range('E29').Value = '  +12.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.44'
$ws.Range('E30').Value = '  +4.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.74'
$ws.Range('E32').Value = '  +7.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.57'
$ws.Range('E33').Value = '  +7.92%  '
$ws.Range('E34').Value = '  +7.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.58'
$ws.Range('E35').Value = '  +9.22%  '
$ws.Range('E36').Value = '  +10.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.96'
$ws.Range('E37').Value = '  +10.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '314.34'
$ws.Range('E38').Value = '  +24.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.869'
$ws.Range('E39').Value = '  +6.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.78'
$ws.Range('E40').Value = '  +11.58%  '
$ws.Range('E41').Value = '  +9.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '35.34'
$ws.Range('E42').Value = '  +4.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.638'
$ws.Range('E43').Value = '  +10.36%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0575'
$ws.Range('E45').Value = '  +10.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.790'
$ws.Range('E46').Value = '  +26.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.989'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.76'
$ws.Range('E48').Value = '  +18.16%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.99'
$ws.Range('E49').Value = '  +13.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0237'
$ws.Range('E50').Value = '  +7.79%  '
$ws.Range('D51').Value = '2.006.02'
$ws.Range('E51').Value = '  +11.51%  '
